# Updated cryptos list data (prices + 1h volume %) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "245.36"); force
# a Text number format before assigning so Excel keeps it as a string
# cell (matching the source inlineStr cells) instead of silently
# converting it to a number, then restore the default "Normal" style so
# no stray number-format style index is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '36.612.93'
$ws.Range("E2").Value = '  -0.70%  '
Set-TextValue $ws.Range("D3") '2.047.96'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '245.36'
$ws.Range("E5").Value = '  +0.09%  '
Set-TextValue $ws.Range("D6") '0.666'
$ws.Range("E6").Value = '  +1.91%  '
$ws.Range("E7").Value = '  +0.05%  '
Set-TextValue $ws.Range("D8") '55.70'
$ws.Range("E8").Value = '  -3.13%  '
Set-TextValue $ws.Range("D9") '63.67'
$ws.Range("E9").Value = '  +8.10%  '
Set-TextValue $ws.Range("D10") '0.368'
$ws.Range("E10").Value = '  -0.28%  '
Set-TextValue $ws.Range("D11") '0.0749'
$ws.Range("E11").Value = '  -3.64%  '
$ws.Range("E12").Value = '  -3.24%  '
Set-TextValue $ws.Range("D13") '0.952'
$ws.Range("E13").Value = '  +8.96%  '
Set-TextValue $ws.Range("D14") '14.74'
$ws.Range("E14").Value = '  -2.36%  '
Set-TextValue $ws.Range("D15") '2.342.20'
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("E16").Value = '  -2.81%  '
Set-TextValue $ws.Range("D17") '2.044.70'
$ws.Range("E17").Value = '  +0.55%  '
Set-TextValue $ws.Range("D18") '36.501.95'
$ws.Range("E18").Value = '  -0.90%  '
Set-TextValue $ws.Range("D19") '17.39'
$ws.Range("E19").Value = '  -0.14%  '
Set-TextValue $ws.Range("D20") '71.99'
$ws.Range("E20").Value = '  -1.51%  '
Set-TextValue $ws.Range("D21") '0.0₃0857'
$ws.Range("E21").Value = '  -3.56%  '
Set-TextValue $ws.Range("D22") '238.06'
$ws.Range("E22").Value = '  +0.86%  '
Set-TextValue $ws.Range("D23") '5.19'
$ws.Range("E23").Value = '  -4.10%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E25").Value = '  -3.43%  '
Set-TextValue $ws.Range("D26") '2.24'
$ws.Range("E26").Value = '  +1.88%  '
$ws.Range("E27").Value = '  -7.74%  '
Set-TextValue $ws.Range("D28") '164.47'
$ws.Range("E28").Value = '  -2.54%  '
Set-TextValue $ws.Range("D29") '20.00'
$ws.Range("E29").Value = '  -0.54%  '
Set-TextValue $ws.Range("D30") '0.122'
$ws.Range("E30").Value = '  -1.94%  '
Set-TextValue $ws.Range("D31") '1.19'
$ws.Range("E31").Value = '  +4.42%  '
$ws.Range("E32").Value = '  -8.43%  '
Set-TextValue $ws.Range("D33") '0.0600'
$ws.Range("E33").Value = '  -2.37%  '
$ws.Range("E34").Value = '  -7.49%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D35") '0.0872'
$ws.Range("E35").Value = '  +2.87%  '
$ws.Range("B36").Value = 'BinanceUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D36") '1.00'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("E38").Value = '  -6.53%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D39") '1.24'
$ws.Range("E39").Value = '  -5.07%  '
$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D40") '4.98'
$ws.Range("E40").Value = '  +2.11%  '
$ws.Range("E41").Value = '  -3.30%  '
Set-TextValue $ws.Range("D42") '2.85'
$ws.Range("E42").Value = '  -3.58%  '
Set-TextValue $ws.Range("D43") '1.11'
$ws.Range("E43").Value = '  -3.16%  '
Set-TextValue $ws.Range("D44") '93.94'
$ws.Range("E44").Value = '  -2.79%  '
Set-TextValue $ws.Range("D45") '0.0902'
$ws.Range("E45").Value = '  -5.36%  '
Set-TextValue $ws.Range("D46") '16.21'
$ws.Range("E46").Value = '  -3.25%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D47") '7.53'
$ws.Range("E47").Value = '  +11.48%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D48") '1.376.44'
$ws.Range("E48").Value = '  +5.39%  '
Set-TextValue $ws.Range("D49") '2.94'
$ws.Range("E49").Value = '  +3.09%  '
$ws.Range("E50").Value = '  -4.27%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range("D51") '45.82'
$ws.Range("E51").Value = '  +1.13%  '
